$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.879.89'
$ws.Range("E2").Value = '  -3.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.181.43'
$ws.Range("E3").Value = '  -3.57%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.03'
$ws.Range("E5").Value = '  -2.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.96'
$ws.Range("E6").Value = '  -6.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.611'
$ws.Range("E7").Value = '  -6.35%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.185.00'
$ws.Range("E9").Value = '  -3.37%  '

$ws.Range("E10").Value = '  -3.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.73'
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  -5.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.739.68'
$ws.Range("E13").Value = '  -3.35%  '

$ws.Range("E14").Value = '  -1.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.010.96'
$ws.Range("E15").Value = '  -3.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.37'
$ws.Range("E16").Value = '  -4.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000159'
$ws.Range("E17").Value = '  -2.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.188.41'
$ws.Range("E18").Value = '  -3.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '414.92'
$ws.Range("E19").Value = '  -4.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.36'
$ws.Range("E20").Value = '  -2.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.79'
$ws.Range("E21").Value = '  -3.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.16'
$ws.Range("E22").Value = '  -3.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.65'
$ws.Range("E24").Value = '  -2.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.202'
$ws.Range("E25").Value = '  +2.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.491'
$ws.Range("E26").Value = '  -3.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000109'
$ws.Range("E27").Value = '  -3.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.75'
$ws.Range("E28").Value = '  -1.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("E30").Value = '  -6.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.87'
$ws.Range("E31").Value = '  -2.32%  '

$ws.Range("E32").Value = '  +0.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.00'
$ws.Range("E33").Value = '  -4.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.39'
$ws.Range("E34").Value = '  -3.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.13'
$ws.Range("E35").Value = '  -5.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.11'
$ws.Range("E36").Value = '  -1.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.37'
$ws.Range("E37").Value = '  -3.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.728.15'
$ws.Range("E38").Value = '  -2.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.70'
$ws.Range("E39").Value = '  -4.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.10'
$ws.Range("E40").Value = '  -5.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.18'
$ws.Range("E41").Value = '  -4.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.719'
$ws.Range("E42").Value = '  -7.11%  '

$ws.Range("E43").Value = '  -3.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.74'
$ws.Range("E44").Value = '  -5.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0628'
$ws.Range("E45").Value = '  -4.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.06'
$ws.Range("E46").Value = '  -5.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '297.04'
$ws.Range("E47").Value = '  -6.86%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0262'
$ws.Range("E48").Value = '  -2.85%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.10'
$ws.Range("E49").Value = '  -9.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0995'
$ws.Range("E50").Value = '  -5.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.10%  '

